$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B1').Value = 'Mercado Libre'
$ws.Range('C1').Value = 'Amazon'
$ws.Range('D1').Value = 'PedidoYa'
$ws.Range('E1').Value = 'Website from Image'
$ws.Range('F1').Value = 'Conclusion'
$ws.Range('B2').Value = 'Utiliza una tipografía sans-serif personalizada, con líneas suavizadas y redondeadas, para una legibilidad óptima y moderna en su plataforma. [8, 40]'
$ws.Range('C2').Value = 'Emplea principalmente Amazon Ember, una sans-serif personalizada, que asegura una lectura clara y profesional en todos sus diversos dispositivos y plataformas. [5, 9, 14]'
$ws.Range('D2').Value = 'Muestra una tipografía sans-serif moderna y limpia en su interfaz, priorizando la claridad y facilidad de lectura para una experiencia de usuario eficiente. [33]'
$ws.Range('E2').Value = 'Utiliza una tipografía sans-serif limpia y estándar, asegurando excelente legibilidad para títulos, descripciones y elementos de navegación de la página.'
$ws.Range('F2').Value = 'Explorar variaciones de fuentes que refuercen la identidad de la marca, manteniendo la claridad y adaptándose mejor a la estética tecnológica.'
$ws.Range('B3').Value = 'Su paleta distintiva combina amarillo brillante, azul profundo y blanco, simbolizando energía, confianza y una plataforma abierta para todos. [2, 3, 8, 19]'
$ws.Range('C3').Value = 'Predominan el naranja energético, el negro y el blanco; estos colores clave crean una identidad de marca reconocible y funcional en su interfaz. [1, 4, 6, 15]'
$ws.Range('D3').Value = 'Ha evolucionado su paleta cromática con un enfoque moderno, usando colores vibrantes y frescos que reflejan su estrategia expansiva y espontánea. [11, 33]'
$ws.Range('E3').Value = 'La paleta principal es azul vibrante para el encabezado, blanco de fondo y negro para el texto, ofreciendo un contraste nítido y funcional.'
$ws.Range('F3').Value = 'Introducir sutilmente colores secundarios para destacar llamadas a la acción o elementos interactivos podría mejorar la jerarquía visual del diseño.'
$ws.Range('B4').Value = 'Mantiene un equilibrio, con un diseño que transmite profesionalismo en transacciones y una accesibilidad amigable que fomenta la comunidad. [2, 8, 32]'
$ws.Range('C4').Value = 'Su diseño es profesional y metódico, centrado en la eficiencia de la compra, aunque mantiene un tono accesible para el cliente. [5, 18, 31]'
$ws.Range('D4').Value = 'Ha adoptado una estrategia más espontánea y abierta, inclinándose hacia un estilo informal que promueve cercanía y facilidad de uso. [33]'
$ws.Range('E4').Value = 'El diseño es profesional y directo, centrado en la funcionalidad de comparación y armado de componentes de PC, transmitiendo seriedad.'
$ws.Range('F4').Value = 'Adoptar un tono más amigable en descripciones o mensajes podría conectar mejor con usuarios menos técnicos, manteniendo profesionalismo.'
$ws.Range('B5').Value = 'Presenta el icónico logotipo de dos manos estrechándose, que simboliza confianza, alianza y una sólida conexión dentro de su comunidad. [2, 8]'
$ws.Range('C5').Value = 'Destaca por la flecha de la sonrisa que va de la ''A'' a la ''Z'' en su logo, sugiriendo variedad y la satisfacción completa del cliente. [5]'
$ws.Range('D5').Value = 'Utiliza un logotipo moderno que incluye elementos visuales distintivos, reforzando su identidad como una plataforma ágil de delivery. [45]'
$ws.Range('E5').Value = 'Un logo de ratón estilizado con ''Rity'', e iconos de búsqueda y perfil genéricos, son claros y funcionales.'
$ws.Range('F5').Value = 'El diseño de iconos más únicos o un emblema secundario podría reforzar la singularidad de la marca en el nicho tecnológico.'
$ws.Range('B6').Value = 'Integra funcionalidades para asegurar la usabilidad general, aunque detalles específicos de accesibilidad no fueron encontrados en la búsqueda. [2]'
$ws.Range('C6').Value = 'Ofrece extensas funciones como VoiceView, textos ajustables y descripciones de audio, garantizando que su plataforma sea inclusiva para todos los usuarios. [7, 12, 25, 31]'
$ws.Range('D6').Value = 'Su rebranding abordó problemas previos de accesibilidad y contraste en la paleta de colores, buscando una mayor legibilidad para los usuarios. [33]'
$ws.Range('E6').Value = 'La combinación de colores y tamaños de fuente ofrece buena legibilidad. No se aprecian características avanzadas de accesibilidad en la interfaz.'
$ws.Range('F6').Value = 'Implementar opciones de contraste alto o herramientas de lectura de pantalla mejoraría la experiencia para usuarios con diversas necesidades visuales.'
$ws.Range('B7').Value = 'Cuenta con una interfaz de usuario intuitiva que guía a los usuarios, con botones clave para búsquedas y gestión de compras, facilitando el recorrido. [2]'
$ws.Range('C7').Value = 'Tiene una navegación fluida con una barra de búsqueda central, botones de compra rápida y un menú claro para explorar amplias categorías de productos. [18, 21, 31]'
$ws.Range('D7').Value = 'Ofrece una navegación sencilla, con botones destacados para pedidos, seguimiento y opciones de cuenta, optimizando la experiencia de entrega rápida.'
$ws.Range('E7').Value = 'La barra superior presenta enlaces vitales: ''Información'', ''Arma tu PC'', ''Comparar'', y una barra de búsqueda prominente.'
$ws.Range('F7').Value = 'Una disposición mejorada de los elementos clave, destacando ''Arma tu PC'', podría optimizar el flujo de usuarios.'
$ws.Range('B8').Value = 'Opera con una estructura jerárquica y funcional, con divisiones claras en tecnología, operaciones y finanzas, optimizando su gestión regional. [27, 29, 43]'
$ws.Range('C8').Value = 'Presenta una estructura organizacional funcional y jerárquica, con equipos especializados en funciones y divisiones geográficas, permitiendo una gestión global eficiente. [23, 24, 35]'
$ws.Range('D8').Value = 'Como parte de Delivery Hero, su organización interna se enfoca en la eficiencia operativa para gestionar entregas en múltiples mercados. [47]'
$ws.Range('E8').Value = 'La página está organizada con un encabezado fijo, título principal y componentes populares en un diseño de cuadrícula estructurada.'
$ws.Range('F8').Value = 'Añadir filtros o categorías claras para los componentes populares facilitaría la exploración y la búsqueda de productos específicos.'
$ws.Range('B9').Value = 'Ofrece un ecosistema integral de servicios como Mercado Pago y Mercado Envíos, expandiendo sus capacidades más allá de las compras básicas. [2]'
$ws.Range('C9').Value = 'Destaca por sus recomendaciones personalizadas, pedidos con un clic, y características innovadoras como Amazon Go y el asistente virtual Alexa. [20, 21, 31]'
$ws.Range('D9').Value = 'Permite el seguimiento de pedidos en tiempo real, opciones de pago variadas y personalización de pedidos, enriqueciendo la experiencia de delivery.'
$ws.Range('E9').Value = 'Las funciones de ''Arma tu PC'' y ''Comparar'' indican herramientas interactivas avanzadas para configurar y evaluar hardware.'
$ws.Range('F9').Value = 'Destacar visiblemente las características únicas como ''Arma tu PC'' con descripciones claras invitaría a más usuarios a explorarlas.'
$ws.Range('B10').Value = 'La plataforma integra guías claras para vendedores sobre cómo usar sus herramientas, lo que facilita la creación y gestión de tiendas online. [44]'
$ws.Range('C10').Value = 'Proporciona una bienvenida opcional para nuevos usuarios y asistentes de voz que incluyen tutoriales en pantalla, facilitando la adaptación del usuario. [20, 28]'
$ws.Range('D10').Value = 'Incorpora instrucciones sencillas para realizar pedidos, personalizar entregas y utilizar las funcionalidades de la aplicación de manera intuitiva.'
$ws.Range('E10').Value = 'El título ''Aprendé, Armá y Compará'' sugiere una sección educativa, aunque no hay tutoriales explícitos visibles aquí.'
$ws.Range('F10').Value = 'Ofrecer mini-guías o tooltips interactivos sobre cómo usar las funciones de armado y comparación podría enriquecer la experiencia.'
$ws.Range('A11').Value = 'User Interface Cohesion'
$ws.Range('B11').Value = 'A pesar de una arquitectura inicial fragmentada, ha logrado una identidad unificada con elementos visuales consistentes en todo su ecosistema. [2, 44]'
$ws.Range('C11').Value = 'Posee una cohesión de interfaz excepcional, manteniendo un diseño limpio y consistente en el sitio web, aplicaciones y materiales de marketing globalmente. [5, 18, 31]'
$ws.Range('D11').Value = 'Su reciente proceso de rebranding buscó una identidad visual más holística, unificando la experiencia en todas sus plataformas digitales. [33]'
$ws.Range('E11').Value = 'El diseño demuestra consistencia en tipografía, colores y disposición de elementos, generando una experiencia visual unificada y ordenada.'
$ws.Range('F11').Value = 'Reforzar la identidad visual con elementos gráficos recurrentes, quizás en los cuadros de producto, podría aumentar la sensación de marca.'
